$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered) from H1 into the new I1/J1 headers
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)
$ws.Cells.Item(1, 10).PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data for columns I (I0) and J (IF), rows 2-71
$iValues = @(
    6,
    5,
    2,
    9,
    9,
    8,
    6,
    7,
    7,
    7,
    8,
    6,
    5,
    5,
    5,
    7,
    3,
    8,
    7,
    6,
    11,
    6,
    7,
    6,
    7,
    9,
    6,
    5,
    7,
    7,
    6,
    6,
    2,
    5,
    8,
    6,
    9,
    6,
    10,
    6,
    7,
    7,
    6,
    7,
    2,
    7,
    1,
    9,
    7,
    7,
    5,
    7,
    8,
    9,
    7,
    8,
    4,
    9,
    6,
    10,
    9,
    8,
    8,
    8,
    5,
    1,
    6,
    6,
    6,
    4
)

$jValues = @(
    7,
    6,
    3,
    9,
    9,
    8,
    6,
    7,
    7,
    7,
    8,
    6,
    6,
    5,
    5,
    7,
    3,
    8,
    7,
    7,
    11,
    6,
    7,
    6,
    7,
    10,
    6,
    5,
    7,
    7,
    6,
    6,
    4,
    5,
    8,
    6,
    9,
    7,
    10,
    7,
    7,
    7,
    7,
    7,
    3,
    7,
    1,
    9,
    9,
    7,
    6,
    8,
    8,
    9,
    7,
    8,
    4,
    9,
    6,
    10,
    9,
    8,
    8,
    8,
    5,
    1,
    6,
    6,
    6,
    4
)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $r = $k + 2
    $ws.Cells.Item($r, 9).Value = $iValues[$k]
    $ws.Cells.Item($r, 10).Value = $jValues[$k]
}

$wb.Save()